# RMA Complete Flow (Repair) - SO To RMA Receipt TO Replacement SO
# Adds a new RMA test-case run (RMA-NNPL-*) replacing the previous
# RMA-QT1Q-* block on the "RMA Details Maintenance Grid" sheet, and
# updates the active sheet/selection state as left behind by the run.

$wb = $excel.ActiveWorkbook

$wsHeader  = $wb.Worksheets.Item("RMA Header")
$wsGrid    = $wb.Worksheets.Item("RMA Details Maintenance Grid")
$wsReceipt = $wb.Worksheets.Item("Receipt")

# ---------------------------------------------------------------
# 1. Update the RMA Details Maintenance Grid test data (row 2-4)
#    Old run: RMA-QT1Q-001/002/003 -> New run: RMA-NNPL-001/002/003
# ---------------------------------------------------------------
$wsGrid.Range("E2").Value = "RMA-NNPL-001"
$wsGrid.Range("F2").Value = "RMA-NNPL-1-1"
$wsGrid.Range("J2").Value = "a7s5f000000xK59AAE"

$wsGrid.Range("E3").Value = "RMA-NNPL-002"
$wsGrid.Range("F3").Value = "RMA-NNPL-1-2"
$wsGrid.Range("J3").Value = "a7s5f000000xK5AAAU"

$wsGrid.Range("E4").Value = "RMA-NNPL-003"
$wsGrid.Range("F4").Value = "RMA-NNPL-1-3"
$wsGrid.Range("J4").Value = "a7s5f000000xK5BAAU"

# ---------------------------------------------------------------
# 2. Restore the sheet-view / selection state left over from the run:
#    - RMA Header:  selection moves to E5
#    - RMA Details Maintenance Grid: becomes the active tab, selection C5
#    - Receipt: no longer the active tab (selection unchanged at A5)
# ---------------------------------------------------------------
$wsHeader.Activate()
$wsHeader.Range("E5").Select()

$wsGrid.Activate()
$wsGrid.Range("C5").Select()
